$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right before "总计", modeled
#    on the existing "2021-Q4" sheet (same column layout/styles).
# ---------------------------------------------------------------
$src       = $wb.Worksheets.Item("2021-Q4")
$totalsRef = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($totalsRef)
$newSheet.Name = "2022-Q1"

# NOTE: once the new sheet has been inserted at the old position of
# "总计", re-resolve the totals sheet by name -- the handle obtained
# before the insert now tracks the newly added sheet instead.
$totals = $wb.Worksheets.Item("总计")

# Copy header row formatting (B1:H1) and the style of A2 from the
# 2021-Q4 sheet so the new sheet matches the existing look & feel.
$src.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)

$src.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Header row text
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row -- B2/D2/E2/F2/G2 are text in the source data (leading
# zeros / fixed decimals), so force text format before assigning,
# then clear the format again so no stray style index is left behind.
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "004641"
$newSheet.Range("C2").Value = "万家量化睿选灵活配置混合"
$newSheet.Range("D2").Value = "0.16"
$newSheet.Range("E2").Value = "85.90"
$newSheet.Range("F2").Value = "1.35"
$newSheet.Range("G2").Value = "0.0022"
$newSheet.Range("H2").Value = 8
$newSheet.Range("B2:G2").ClearFormats()

# ---------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: a new first data row for
#    2022-Q1 is inserted above the existing history, which shifts
#    down by one row.
# ---------------------------------------------------------------
$totals.Rows("2:2").Insert()

# The inserted row inherits formatting from the row above (the bold
# header) -- strip it back to plain, then re-apply the same styling
# the other index cells (A3/A4 = style of A1's column) carry.
$totals.Range("A2:D2").ClearFormats()
$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 1
$totals.Range("D2").Value = 0

$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2

Write-Host "done"
